# Narea_model.xlsx — "fixed slopes for log transformed variables"
#
# Corrects the regression Slope/SE values for three log-transformed
# predictors (ln PAR, ln VPD, ln LMA), restyles the RelImp (column E)
# values to match the p-value number format, bumps the Elevation
# row's Slope/SE to a 4-decimal format, fills in the three trailing
# RelImp "NA" placeholders with their real values, and removes the
# now-empty trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected regression coefficients (Slope / SE) -----------------
# NB: the PowerShell parser here does not accept scientific-notation
# numeric literals (e.g. "1.2E-3"), so every value below is written out
# in plain decimal form.

# ln PAR (row 7)
$ws.Range("B7").Value = -0.114992084861493
$ws.Range("C7").Value = 0.285740918369213

# ln VPD (row 8)
$ws.Range("B8").Value = -0.033297323192085997
$ws.Range("C8").Value = 0.085705776179994406

# ln LMA (row 10)
$ws.Range("B10").Value = 0.93553292553719403
$ws.Range("C10").Value = 0.0087712118158288204

# --- Fill in the previously-missing RelImp values --------------------
$ws.Range("E13").Value = 0.13414764901703899
$ws.Range("E14").Value = 0.149452310426681
$ws.Range("E15").Value = 0.048506400071264201
$ws.Range("E16").Value = 0.069402583381629704

# --- Number formats ----------------------------------------------------
# RelImp data (E2:E16) now matches the "p" column's 3-decimal format.
$ws.Range("E2:E16").NumberFormat = "0.000"
# RelImp header (E1) gets its own distinct 1-decimal format.
$ws.Range("E1").NumberFormat = "0.0"
# Elevation row Slope/SE (B9:C9) gets a 4-decimal format.
$ws.Range("B9:C9").NumberFormat = "0.0000"

# --- Column widths -------------------------------------------------------
# B and C now share one (slightly wider) best-fit width.
$ws.Range("B1:C16").ColumnWidth = 5.83

# --- Remove the now-empty trailing row ------------------------------------
$ws.Rows(17).Delete()

# --- Selection cosmetic state ---------------------------------------------
$ws.Range("F21").Select()
